$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A2").Value = "Diffa"
$ws.Range("B2").Value = "P"
$ws.Range("D2").Value = 2123434565
$ws.Range("A3").Value = "Diffa"
$ws.Range("D3").Value = 2123434565
$ws.Range("E3").Value = 12345678
$ws.Range("F3").Value = 12345678
[void]$ws.Range("I7").Select()
